# Fix latency units in report sheets
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update header O2: "Utility" -> "Utility (Percent)"
$ws.Range("O2").Value = "Utility (Percent)"

# Append " msec" to Read Latency columns (I, J, K) for data rows 3-23
for ($row = 3; $row -le 23; $row++) {
    foreach ($col in @("I", "J", "K")) {
        $cell = $ws.Range("$col$row")
        $current = $cell.Value2
        $cell.Value = "$current msec"
    }
}
